$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44327

$ws.Range("D3").Value = 44280
$ws.Range("K3").Value = 1800
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = 1900
$ws.Range("P3").Value = 633

$ws.Range("D4").Value = 44280
$ws.Range("I4").Value = "Segunda"
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 1400
$ws.Range("L4").Value = 1500
$ws.Range("M4").Value = 1450
$ws.Range("P4").Value = 483

$ws.Range("D5").Value = 44270
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 1800
$ws.Range("L5").Value = 2000
$ws.Range("M5").Value = 1900
$ws.Range("P5").Value = 633

$ws.Range("D6").Value = 44270
$ws.Range("I6").Value = "Segunda"
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 1200
$ws.Range("L6").Value = 1500
$ws.Range("M6").Value = 1350
$ws.Range("P6").Value = 450

$ws.Range("D7").Value = 44364
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 1700
$ws.Range("L7").Value = 1800
$ws.Range("M7").Value = 1750
$ws.Range("P7").Value = 583

$ws.Range("D8").Value = 44364
$ws.Range("I8").Value = "Segunda"
$ws.Range("K8").Value = 1400
$ws.Range("L8").Value = 1500
$ws.Range("M8").Value = 1450
$ws.Range("P8").Value = 483

$ws.Range("D9").Value = 44306
$ws.Range("I9").Value = "Primera"
$ws.Range("K9").Value = 2400
$ws.Range("L9").Value = 2500
$ws.Range("M9").Value = 2450
$ws.Range("P9").Value = 817

$ws.Range("D10").Value = 44333
$ws.Range("L10").Value = 1700
$ws.Range("M10").Value = 1600
$ws.Range("P10").Value = 533

$ws.Range("D11").Value = 44166
$ws.Range("J11").Value = 250
$ws.Range("K11").Value = 900
$ws.Range("L11").Value = 1000
$ws.Range("M11").Value = 950
$ws.Range("P11").Value = 317

$ws.Range("D12").Value = 44431
$ws.Range("J12").Value = 250
$ws.Range("K12").Value = 1000
$ws.Range("L12").Value = 1300
$ws.Range("M12").Value = 1150
$ws.Range("P12").Value = 383

$ws.Range("D13").Value = 44174
$ws.Range("J13").Value = 250
$ws.Range("K13").Value = 500
$ws.Range("L13").Value = 600
$ws.Range("M13").Value = 550
$ws.Range("P13").Value = 183

$ws.Range("D14").Value = 44278
$ws.Range("J14").Value = 140
$ws.Range("K14").Value = 2000
$ws.Range("L14").Value = 2500
$ws.Range("M14").Value = 2250
$ws.Range("P14").Value = 750

$ws.Range("D15").Value = 44278
$ws.Range("I15").Value = "Segunda"
$ws.Range("K15").Value = 1500
$ws.Range("M15").Value = 1650
$ws.Range("P15").Value = 550

$ws.Range("D16").Value = 44385
$ws.Range("I16").Value = "Primera"
$ws.Range("K16").Value = 2000
$ws.Range("L16").Value = 2300
$ws.Range("M16").Value = 2150
$ws.Range("P16").Value = 717

$ws.Range("D17").Value = 44300
$ws.Range("J17").Value = 250
$ws.Range("K17").Value = 1600
$ws.Range("L17").Value = 1800
$ws.Range("M17").Value = 1700
$ws.Range("P17").Value = 567

$ws.Range("D18").Value = 44224
$ws.Range("I18").Value = "Primera"
$ws.Range("K18").Value = 1400
$ws.Range("L18").Value = 1500
$ws.Range("M18").Value = 1450
$ws.Range("P18").Value = 483

$ws.Range("D19").Value = 44224
$ws.Range("I19").Value = "Segunda"
$ws.Range("J19").Value = 160
$ws.Range("L19").Value = 1200
$ws.Range("M19").Value = 1100
$ws.Range("P19").Value = 367

$ws.Range("D20").Value = 44295
$ws.Range("K20").Value = 1500
$ws.Range("L20").Value = 1800
$ws.Range("M20").Value = 1650
$ws.Range("P20").Value = 550

$ws.Range("D21").Value = 44249
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 100
$ws.Range("K21").Value = 1500
$ws.Range("L21").Value = 1600
$ws.Range("M21").Value = 1550
$ws.Range("P21").Value = 517

$ws.Range("D22").Value = 44397

$ws.Range("D23").Value = 44398
$ws.Range("J23").Value = 300
$ws.Range("K23").Value = 1700
$ws.Range("L23").Value = 1800
$ws.Range("M23").Value = 1750
$ws.Range("P23").Value = 583

$ws.Range("D24").Value = 44417
$ws.Range("J24").Value = 250
$ws.Range("K24").Value = 1800
$ws.Range("L24").Value = 2000
$ws.Range("M24").Value = 1900
$ws.Range("P24").Value = 633

$ws.Range("D25").Value = 44417
$ws.Range("L25").Value = 1600
$ws.Range("M25").Value = 1550
$ws.Range("P25").Value = 517

$ws.Range("D26").Value = 44432
$ws.Range("J26").Value = 200
$ws.Range("K26").Value = 1200
$ws.Range("L26").Value = 1300
$ws.Range("M26").Value = 1250
$ws.Range("P26").Value = 417

$ws.Range("D27").Value = 44432
$ws.Range("K27").Value = 950
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = 975
$ws.Range("P27").Value = 325

$ws.Range("D28").Value = 44342
$ws.Range("K28").Value = 2000
$ws.Range("L28").Value = 2200
$ws.Range("M28").Value = 2100
$ws.Range("P28").Value = 700

$ws.Range("D29").Value = 44302
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 200
$ws.Range("K29").Value = 1400
$ws.Range("L29").Value = 1500
$ws.Range("M29").Value = 1450
$ws.Range("P29").Value = 483

$ws.Range("D30").Value = 44391
$ws.Range("J30").Value = 250
$ws.Range("K30").Value = 1800
$ws.Range("L30").Value = 2000
$ws.Range("M30").Value = 1900
$ws.Range("P30").Value = 633

$ws.Range("D31").Value = 44161
$ws.Range("J31").Value = 200
$ws.Range("K31").Value = 600
$ws.Range("L31").Value = 700
$ws.Range("M31").Value = 650
$ws.Range("P31").Value = 217

$ws.Range("D32").Value = 44161
$ws.Range("I32").Value = "Segunda"
$ws.Range("J32").Value = 250
$ws.Range("K32").Value = 500
$ws.Range("L32").Value = 600
$ws.Range("M32").Value = 550
$ws.Range("P32").Value = 183

$ws.Range("D33").Value = 44428
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 200
$ws.Range("K33").Value = 1500
$ws.Range("L33").Value = 1800
$ws.Range("M33").Value = 1650
$ws.Range("P33").Value = 550

$ws.Range("D34").Value = 44435
$ws.Range("J34").Value = 450
$ws.Range("K34").Value = 1000
$ws.Range("L34").Value = 1300
$ws.Range("M34").Value = 1194
$ws.Range("P34").Value = 398

$ws.Range("D35").Value = 44435
$ws.Range("I35").Value = "Segunda"
$ws.Range("K35").Value = 950
$ws.Range("L35").Value = 1000
$ws.Range("M35").Value = 975
$ws.Range("P35").Value = 325

$ws.Range("D36").Value = 44376
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 340
$ws.Range("K36").Value = 1400
$ws.Range("L36").Value = 1500
$ws.Range("M36").Value = 1471
$ws.Range("P36").Value = 490
